$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 40
$ws.Range("I6").Value = 20
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 60
$ws.Range("L6").Value = 300
$ws.Range("M6").Value = 52
$ws.Range("N6").Value = -524
$ws.Range("H17").Value = 1888.25
$ws.Range("J17").Value = 2025.9
$ws.Range("L17").Value = 6077.700000000001
$ws.Range("N17").Value = -6413.700000000001
$ws.Range("H32").Value = 19891.715
$ws.Range("I32").Value = 13249
$ws.Range("J32").Value = 22548.8
$ws.Range("K32").Value = 13249
$ws.Range("L32").Value = 22548.8
$ws.Range("M32").Value = -12923
$ws.Range("N32").Value = -23200.8
$ws.Range("H58").Value = 2092.2
$ws.Range("I58").Value = 303.42856
$ws.Range("J58").Value = 6266
$ws.Range("K58").Value = 910.28568
$ws.Range("L58").Value = 18798
$ws.Range("M58").Value = -760.28568
$ws.Range("N58").Value = -19098
$ws.Range("H62").Value = 20842396
$ws.Range("J62").Value = 3874.75
$ws.Range("L62").Value = 3874.75
$ws.Range("N62").Value = -5122.75
$ws.Range("H65").Value = 20842396
$ws.Range("J65").Value = 3874.75
$ws.Range("L65").Value = 19373.75
$ws.Range("N65").Value = -25613.75
$ws.Range("H98").Value = 3373.7693
$ws.Range("I98").Value = 2805.818
$ws.Range("K98").Value = 2805.818
$ws.Range("M98").Value = -1307.818
$ws.Range("H122").Value = 3373.7693
$ws.Range("I122").Value = 2805.818
$ws.Range("K122").Value = 8417.454000000002
$ws.Range("M122").Value = -5967.454000000002
$ws.Range("H132").Value = 265352
$ws.Range("I132").Value = 287241.78
$ws.Range("K132").Value = 861725.3400000001
$ws.Range("M132").Value = -859195.3400000001
$ws.Range("H137").Value = 8916.4
$ws.Range("I137").Value = 11324
$ws.Range("K137").Value = 33972
$ws.Range("M137").Value = -31422

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4733.9546
$ws.Range("I2").Value = 2074
$ws.Range("J2").Value = 11827.167
$ws.Range("K2").Value = 2074
$ws.Range("L2").Value = 11827.167
$ws.Range("M2").Value = -1961
$ws.Range("N2").Value = -12053.167
$ws.Range("H61").Value = 4723.4917
$ws.Range("I61").Value = 4746.1055
$ws.Range("J61").Value = 4713.2617
$ws.Range("K61").Value = 4746.1055
$ws.Range("L61").Value = 4713.2617
$ws.Range("M61").Value = -4534.1055
$ws.Range("N61").Value = -5137.2617
$ws.Range("H80").Value = 10555.5
$ws.Range("J80").Value = 11111
$ws.Range("L80").Value = 11111
$ws.Range("N80").Value = -13107
$ws.Range("H83").Value = 10555.5
$ws.Range("J83").Value = 11111
$ws.Range("L83").Value = 33333
$ws.Range("N83").Value = -43317
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H116").Value = 4733.9546
$ws.Range("I116").Value = 2074
$ws.Range("J116").Value = 11827.167
$ws.Range("K116").Value = 2074
$ws.Range("L116").Value = 11827.167
$ws.Range("M116").Value = 220
$ws.Range("N116").Value = -16415.167
$ws.Range("H122").Value = 2959.484
$ws.Range("I122").Value = 2561.318
$ws.Range("K122").Value = 7683.954000000001
$ws.Range("M122").Value = -5233.954000000001
$ws.Range("H132").Value = 525575.3
$ws.Range("I132").Value = 527463.5600000001
$ws.Range("J132").Value = 461375
$ws.Range("K132").Value = 1582390.68
$ws.Range("L132").Value = 1384125
$ws.Range("M132").Value = -1579860.68
$ws.Range("N132").Value = -1389185
$ws.Range("H135").Value = 53374.5
$ws.Range("J135").Value = 53374.5
$ws.Range("L135").Value = 53374.5
$ws.Range("N135").Value = -63514.5
$ws.Range("H136").Value = 4723.4917
$ws.Range("I136").Value = 4746.1055
$ws.Range("J136").Value = 4713.2617
$ws.Range("K136").Value = 14238.3165
$ws.Range("L136").Value = 14139.7851
$ws.Range("M136").Value = -11688.3165
$ws.Range("N136").Value = -19239.7851

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4733.9546
$ws.Range("I3").Value = 2074
$ws.Range("J3").Value = 11827.167
$ws.Range("K3").Value = 2074
$ws.Range("L3").Value = 11827.167
$ws.Range("M3").Value = -1960
$ws.Range("N3").Value = -12055.167
$ws.Range("H86").Value = 3950.9565
$ws.Range("I86").Value = 1805.3334
$ws.Range("K86").Value = 1805.3334
$ws.Range("M86").Value = -682.3334
$ws.Range("H89").Value = 3950.9565
$ws.Range("I89").Value = 1805.3334
$ws.Range("K89").Value = 9026.666999999999
$ws.Range("M89").Value = -3410.666999999999
$ws.Range("H94").Value = 5282.409
$ws.Range("I94").Value = 2152.6
$ws.Range("J94").Value = 7890.5835
$ws.Range("K94").Value = 2152.6
$ws.Range("L94").Value = 7890.5835
$ws.Range("M94").Value = -1701.6
$ws.Range("N94").Value = -8792.583500000001
$ws.Range("H105").Value = 2491.75
$ws.Range("I105").Value = 2658.3333
$ws.Range("J105").Value = 1992
$ws.Range("K105").Value = 2658.3333
$ws.Range("L105").Value = 1992
$ws.Range("M105").Value = -911.3332999999998
$ws.Range("N105").Value = -5486
$ws.Range("H132").Value = 95696.75
$ws.Range("J132").Value = 95696.75
$ws.Range("L132").Value = 95696.75
$ws.Range("N132").Value = -105816.75
$ws.Range("H134").Value = 1196323.5
$ws.Range("I134").Value = 1393942
$ws.Range("J134").Value = 10612.833
$ws.Range("K134").Value = 4181826
$ws.Range("L134").Value = 31838.499
$ws.Range("M134").Value = -4179291
$ws.Range("N134").Value = -36908.499

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6954.9287
$ws.Range("I31").Value = 1455.4286
$ws.Range("J31").Value = 12454.429
$ws.Range("K31").Value = 1455.4286
$ws.Range("L31").Value = 12454.429
$ws.Range("M31").Value = -1160.4286
$ws.Range("N31").Value = -13044.429
$ws.Range("H34").Value = 6954.9287
$ws.Range("I34").Value = 1455.4286
$ws.Range("J34").Value = 12454.429
$ws.Range("K34").Value = 1455.4286
$ws.Range("L34").Value = 12454.429
$ws.Range("M34").Value = -1253.4286
$ws.Range("N34").Value = -12858.429
$ws.Range("H132").Value = 16270.889
$ws.Range("I132").Value = 6309.273
$ws.Range("J132").Value = 31924.857
$ws.Range("K132").Value = 18927.819
$ws.Range("L132").Value = 95774.571
$ws.Range("M132").Value = -16397.819
$ws.Range("N132").Value = -100834.571
$ws.Range("H134").Value = 58833644
$ws.Range("I134").Value = 71433800
$ws.Range("K134").Value = 214301400
$ws.Range("M134").Value = -214298865

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 14010
$ws.Range("I3").Value = 2622.8572
$ws.Range("J3").Value = 29952
$ws.Range("K3").Value = 7868.571599999999
$ws.Range("L3").Value = 89856
$ws.Range("M3").Value = -7756.571599999999
$ws.Range("N3").Value = -90080
$ws.Range("H12").Value = 952975.9
$ws.Range("J12").Value = 999.3
$ws.Range("L12").Value = 2997.9
$ws.Range("N12").Value = -3343.9
$ws.Range("H107").Value = 3650.4524
$ws.Range("I107").Value = 957.5
$ws.Range("J107").Value = 4099.278
$ws.Range("K107").Value = 2872.5
$ws.Range("L107").Value = 12297.834
$ws.Range("M107").Value = -952.5
$ws.Range("N107").Value = -16137.834
$ws.Range("H124").Value = 33124.5
$ws.Range("I124").Value = 27499.334
$ws.Range("K124").Value = 82498.00199999999
$ws.Range("M124").Value = -77588.00199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H45").Value = 26774.666
$ws.Range("J45").Value = 30326
$ws.Range("L45").Value = 30326
$ws.Range("N45").Value = -31444
$ws.Range("H70").Value = 7501.125
$ws.Range("I70").Value = 6670.1665
$ws.Range("K70").Value = 6670.1665
$ws.Range("M70").Value = -6400.1665
$ws.Range("H73").Value = 7501.125
$ws.Range("I73").Value = 6670.1665
$ws.Range("K73").Value = 6670.1665
$ws.Range("M73").Value = -5734.1665
$ws.Range("H132").Value = 5134.8
$ws.Range("I132").Value = 5134.8
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 15404.4
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4938.1304
$ws.Range("I122").Value = 4961.5713
$ws.Range("J122").Value = 4901.6665
$ws.Range("K122").Value = 14884.7139
$ws.Range("L122").Value = 14704.9995
$ws.Range("M122").Value = -12434.7139
$ws.Range("N122").Value = -19604.9995
$ws.Range("H132").Value = 2241.3704
$ws.Range("I132").Value = 2241.3704
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6724.111199999999
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 35724196
$ws.Range("I136").Value = 90918800
$ws.Range("K136").Value = 272756400
$ws.Range("M136").Value = -272753850

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3114.111
$ws.Range("I126").Value = 1899.9333
$ws.Range("K126").Value = 5699.7999
$ws.Range("M126").Value = -3229.7999
$ws.Range("H132").Value = 8843.450000000001
$ws.Range("I132").Value = 5157.4546
$ws.Range("J132").Value = 26220.285
$ws.Range("K132").Value = 15472.3638
$ws.Range("L132").Value = 78660.855
$ws.Range("M132").Value = -12942.3638
$ws.Range("N132").Value = -83720.855
$ws.Range("H136").Value = 21754070
$ws.Range("I136").Value = 23820898
$ws.Range("J136").Value = 52375
$ws.Range("K136").Value = 71462694
$ws.Range("L136").Value = 157125
$ws.Range("M136").Value = -71460144
$ws.Range("N136").Value = -162225
